$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 541, shifting the existing rows 541:658 down to 542:659.
$ws.Rows("541:541").Insert()

# Populate the newly inserted row with the new price-survey record
# (weekly refresh added a new data point for this market/product).
$ws.Range("A541").Value = 3
$ws.Range("B541").Value = "Femacal de La Calera"
$ws.Range("C541").Value = "Coquimbo"
$ws.Range("D541").Value = 45275
$ws.Range("E541").Value = 5
$ws.Range("F541").Value = 100112012
$ws.Range("G541").Value = "Espinaca"
$ws.Range("H541").Value = "Sin especificar"
$ws.Range("I541").Value = "Primera"
$ws.Range("J541").Value = 115
$ws.Range("K541").Value = 5000
$ws.Range("L541").Value = 5500
$ws.Range("M541").Value = 5239
$ws.Range("N541").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O541").Value = "Provincia de Quillota"
$ws.Range("P541").Value = 1746
$ws.Range("Q541").Value = 3
$ws.Range("R541").Value = "Hortaliza"
